# OneR, ID3 and C45 data added
#
# The "Clean Data" sheet holds a tabular dataset (German credit data) with
# columns:
#   A Account Status        B Debt History        C Loan Purpose
#   D Credit Amount (num)   E Savings              F Years Employed
#   G Personal Status And Gender   H Age (num)      I Job Status
#   J Decision (good/bad)
#
# For the new OneR / ID3 / C45 classifiers the categorical feature values
# need to be recognisable as literal nominal labels (rather than things
# Excel/consumers might try to reinterpret, e.g. "<100" or "0<=X<200"), so
# every categorical feature value is wrapped in single quotes. The numeric
# columns (Credit Amount, Age) and the Decision/target column (good/bad)
# are left exactly as they were.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 1000 }

# Columns holding categorical text values that must be quoted.
$targetCols = @(1, 2, 3, 5, 6, 7, 9)

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($c in $targetCols) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()
        if ($v -ne $null -and $v -ne "") {
            # Leading "''" -> a literal single quote followed by the text
            # (a single leading quote would be swallowed by Excel as a
            # "force text" prefix marker instead of becoming real content).
            $cell.Value = "''" + $v + "'"
        }
    }
}
